# StatusTableOverview.xlsx — "Updated info on TicketService and S-Lab"
#
# Table columns: A=Name, B=Rebeca Type, C=Where, D=Check in Afra,
#                E=Code Edits, F=Reached States, G=Reached Transitions,
#                H=PlantUML status, I=Comments

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("StatusTableOverview")

# --- Row 4: Scientific Lab ---------------------------------------------
# Status note got more pessimistic; comment expanded with the real cause.
$ws.Range("D4:E4").Style = "Bad"
$ws.Range("I4").Value = "Deadlock in diagram (more like Deadlock in code) Code needs to be fixed."
$ws.Range("D4").Value = "OK, but needs to be fixed"

# --- Row 9: Ticket Service -----------------------------------------------
# Model-checking run finally produced state/transition counts and a review
# comment was added.
$ws.Range("F9:G9").Style = "Neutral"
$ws.Range("F9").Value = 699
$ws.Range("G9").Value = 827
$ws.Range("I9").Value = 'Waiting for review. But in main creates 2 Customers. The UML is more generic. No changed done to the code. How do GPT know that the constructors contains the variable "id" of the costumer?'

# Reflect the cursor position left behind by the edit.
$ws.Range("H13").Select()
